# Updated cryptos list - apply Price (D) and Volume(1h) (E) changes per row
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.637.38"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.96%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.669.07"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.02%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "601.69"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.63%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "156.95"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +4.34%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +3.42%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -1.19%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.84"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -2.72%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.31%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "29.42"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -3.07%  "
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -5.61%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.150.73"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.96%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.405.69"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.04%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.668.05"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.18%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.83"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.50%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.80"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -2.39%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.68"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +1.54%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "351.74"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -3.24%  "
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.13%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "69.82"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.44%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000110"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.88%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.75"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.09%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.64"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -2.27%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -4.10%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -5.40%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.80%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.27%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.16"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -2.82%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "532.00"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -1.61%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -3.07%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.52"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -1.59%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.51"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.46%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.424"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -2.43%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "20.45"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -2.11%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.01%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "158.59"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -2.97%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -3.85%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.04%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "42.72"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.05%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "164.97"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -3.51%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.11"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -2.91%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.31"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -2.04%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -1.67%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "23.09"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -1.24%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -3.17%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0259"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -3.24%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.101"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "20.25"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -1.31%  "
